$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidentiality footer date from 2021-04-28 to 2021-04-29
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-39
$ws.Range("D2").Value = 0.06218365871294289
$ws.Range("E2").Value = -0.0007486150621351939
$ws.Range("D3").Value = 0.0551171222609106
$ws.Range("E3").Value = -0.008053111250785672
$ws.Range("D4").Value = 0.2966096236063067
$ws.Range("E4").Value = 0.01762902401635147
$ws.Range("D5").Value = 0.03744157906571325
$ws.Range("E5").Value = 0.00370391788347546
$ws.Range("D6").Value = 0.03285516155199648
$ws.Range("E6").Value = 0.01460698222317425
$ws.Range("D7").Value = 0.02966464318667732
$ws.Range("E7").Value = 0.0194442619720161
$ws.Range("D8").Value = 0.02805742131936612
$ws.Range("E8").Value = 0.01370539572786766
$ws.Range("D9").Value = 0.02388466369753916
$ws.Range("E9").Value = 0.01261875407933877
$ws.Range("D10").Value = 0.02553887022673997
$ws.Range("E10").Value = 0.01429395008138923
$ws.Range("D11").Value = 0.02327256397789903
$ws.Range("E11").Value = 0.07297297297297267
$ws.Range("D12").Value = 0.02253488313431432
$ws.Range("E12").Value = 0.02697976517611789
$ws.Range("D13").Value = 0.02183909875190767
$ws.Range("E13").Value = 0.01057854844866135
$ws.Range("D14").Value = 0.0213379651116151
$ws.Range("E14").Value = 0.008447488584474749
$ws.Range("D15").Value = 0.02076700403602576
$ws.Range("E15").Value = 0.01760970879859891
$ws.Range("D16").Value = 0.02135907573158853
$ws.Range("E16").Value = 0.01763856154489485
$ws.Range("D17").Value = 0.02141645331408045
$ws.Range("E17").Value = -0.0169341589788955
$ws.Range("D18").Value = 0.01497100212762571
$ws.Range("E18").Value = 0.01145435612634516
$ws.Range("D19").Value = 0.01639916263372631
$ws.Range("E19").Value = 0.04251386321626627
$ws.Range("D20").Value = 0.01527272160386149
$ws.Range("E20").Value = 0.01577175261385788
$ws.Range("D21").Value = 0.016356508355421
$ws.Range("E21").Value = 0.01428325589399426
$ws.Range("D22").Value = 0.01503509180467329
$ws.Range("E22").Value = -0.0250576036866359
$ws.Range("D23").Value = 0.01508424165835504
$ws.Range("E23").Value = 0.01250233252472466
$ws.Range("D24").Value = 0.01441238264525162
$ws.Range("E24").Value = 0.0122739018087854
$ws.Range("D25").Value = 0.0138328690620833
$ws.Range("E25").Value = 0.006847974955977332
$ws.Range("D26").Value = 0.01446293987359827
$ws.Range("E26").Value = -0.003622890078221475
$ws.Range("D27").Value = 0.01273554986182288
$ws.Range("E27").Value = 0.01950883635529022
$ws.Range("D28").Value = 0.01366928882218655
$ws.Range("E28").Value = 0.01596654628397665
$ws.Range("D29").Value = 0.01428528588705255
$ws.Range("E29").Value = -0.005410976552434832
$ws.Range("D30").Value = 0.01295131122391039
$ws.Range("E30").Value = 0.01471178280058849
$ws.Range("D31").Value = 0.01228356607521205
$ws.Range("E31").Value = 0.01487696538108985
$ws.Range("D32").Value = 0.01332924545122976
$ws.Range("E32").Value = -0.009291521486643473
$ws.Range("D33").Value = 0.01251859764424968
$ws.Range("E33").Value = -0.04423401219354006
$ws.Range("D34").Value = 0.006615418742138325
$ws.Range("E34").Value = 0.003142029554715542
$ws.Range("D35").Value = 0.00548356473279314
$ws.Range("E35").Value = 0.004896154149885534
$ws.Range("D36").Value = 0.005869618429333066
$ws.Range("E36").Value = -0.01195175034121487
$ws.Range("D37").Value = 0.005580565325081355
$ws.Range("E37").Value = 0.001183363079071853
$ws.Range("D38").Value = 0.004971280354771006
$ws.Range("E38").Value = 0.01543989547038338
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 0.01160239673740238

$ws.Protect()
